# Final basic episode production
# Updates Performance_SNOW timing data with a new measurement run, adds a
# new "basic episode QA" data row, and refreshes the dependent views.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance_SNOW")

# --- 1. Refresh the raw timestamp columns (B:G) for the existing rows ---
# Each entry is: row, B, C, D, E, F, G  ($null means "leave/make blank")
$data = @(
  @(2,  44265.512569444443, 44265.512754629628, 44266.516689814816, 44266.516956018517, 44266.516956018517, 44266.517141203702),
  @(3,  44265.512754629628, 44265.512777777774, 44266.517141203702, 44266.517175925925, 44266.517175925925, 44266.517175925925),
  @(4,  44265.512777777774, 44265.517592592594, 44266.517175925925, 44266.518159722225, 44266.518159722225, 44266.518252314818),
  @(5,  44265.517592592594, 44265.517685185187, 44266.518252314818, 44266.518333333333, 44266.518333333333, 44266.518391203703),
  @(6,  44265.517685185187, 44265.517789351848, 44266.518391203703, 44266.518449074072, 44266.518449074072, 44266.518553240741),
  @(7,  44265.517789351848, 44265.538356481484, 44266.518553240741, 44266.52270833333,  44266.52270833333,  44266.545011574075),
  @(8,  44213.803611111114, 44213.803738425922, $null,              $null,              44248.301840277774, 44248.302141203705),
  @(9,  44213.804085648146, 44213.804328703707, $null,              $null,              44248.302141203705, 44248.302349537036),
  @(10, 44213.804386574076, 44213.804490740738, $null,              $null,              44248.302349537036, 44248.302499999998),
  @(11, 44213.804664351854, 44213.804837962962, $null,              $null,              44248.302499999998, 44248.302905092591),
  @(12, 44265.538356481484, 44265.547210648147, 44266.545011574075, 44266.546307870369, 44266.546307870369, 44266.552370115744),
  @(13, 44265.547210648147, 44265.547731481478, 44266.552361111113, 44266.552499999998, 44266.552499999998, 44266.552835648145),
  @(14, 44265.547731481478, 44265.563562604169, 44266.552835648145, 44266.555677939818, 44266.555671296293, 44266.569363425922),
  @(15, 44265.56355324074,  44265.75880787037,  44266.569363425922, 44266.601527777777, 44266.601527777777, 44266.767002314817),
  @(16, 44265.75880787037,  44265.974907407406, 44266.767002314817, 44266.807393750001, 44266.807384259257, 44266.97210648148),
  @(17, 44266.332280092596, 44266.418611111112, 44266.97210648148,  44266.988981481481, 44266.988981481481, 44267.047349537039),
  @(18, 44266.418611111112, 44266.46298611111,  44267.047349537039, 44267.056574074071, 44267.056574074071, 44267.085613425923)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -eq $null) {
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -eq $null) {
        $ws.Cells.Item($r, 5).ClearContents()
    } else {
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# --- 2. Insert the new measurement row above the totals row ---
# This pushes the existing totals (previously row 19) down to row 20 and
# leaves a fresh row 19 ready for the new "basic episode QA" record.
$ws.Rows("19").Insert()

$ws.Cells.Item(19, 1).Value = "basic episode QA"
$ws.Cells.Item(19, 2).Value = 44266.46298611111
$ws.Cells.Item(19, 3).Value = 44266.463055555556
$ws.Cells.Item(19, 4).Value = 44267.085613425923
$ws.Cells.Item(19, 5).Value = 44267.085729166669
$ws.Cells.Item(19, 6).Value = 44267.085729166669
$ws.Cells.Item(19, 7).Value = 44267.085787037038
$ws.Cells.Item(19, 13).Value = 0.00010416666918899864

# --- 3. (Re)write the H:K elapsed-time formulas for every data row ---
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 8).Formula  = "=C" + $r + "-B" + $r
    $ws.Cells.Item($r, 9).Formula  = "=E" + $r + "-D" + $r
    $ws.Cells.Item($r, 10).Formula = "=G" + $r + "-F" + $r
    $ws.Cells.Item($r, 11).Formula = "=SUM(H" + $r + ":J" + $r + ")"
}

# --- 4. Totals row now lives at row 20 and must sum through row 19 ---
$ws.Cells.Item(20, 8).Formula  = "=SUM(H2:H19)"
$ws.Cells.Item(20, 9).Formula  = "=SUM(I2:I19)"
$ws.Cells.Item(20, 10).Formula = "=SUM(J2:J19)"
$ws.Cells.Item(20, 11).Formula = "=SUM(H20:J20)"

# --- 5. Sheet-level view/formatting touch-ups ---
$ws.Columns("N").ColumnWidth = 7.6666666666667
$ws.Range("J25").Select()

# Pentaho_CHPC keeps its frozen-pane view but the remembered selection moves
$ws4 = $wb.Worksheets.Item("Pentaho_CHPC")
$ws4.Range("I48").Select()
$ws.Activate()

Write-Host "edit applied"
